# Insert a new data row at row 6 (pushing the existing rows 6-49 down to 7-50),
# and populate the newly inserted row with the new price-record values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 6..49 down to 7..50, creating a blank row 6.
$ws.Rows.Item(6).Insert()

# Fill in the new row 6 with the new record's data.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Macroferia Regional de Talca"
$ws.Range("C6").Value = "Maule"
$ws.Range("D6").Value = 44503
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 100112026
$ws.Range("G6").Value = "Haba"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 6000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 6000
$ws.Range("N6").Value = "`$/saco 25 kilos"
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 240
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
